$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$dst = $ws.Cells.Item(158, 8)
$ws.Range("H5").Copy($ws.Range("H158"))
$dst.Value = 0
$dst.Font.Name = "Calibri"
$dst.Font.Size = 11
